$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 21
$ws.Range("C4").Value = 24
$ws.Range("C5").Value = 32
$ws.Range("C6").Value = 30
$ws.Range("C7").Value = 29
$ws.Range("C8").Value = 31
$ws.Range("C9").Value = 20
$ws.Range("C10").Value = 28
$ws.Range("C11").Value = 33
$ws.Range("C12").Value = 27
$ws.Range("C14").Value = 35
$ws.Range("C15").Value = 7

# Update text values in column B
$ws.Range("B5").Value = "<long>"
$ws.Range("B9").Value = "<nove>"
